# Daily attendance processing - 2026-01-28 14:46:05
# Swap the order of "dnasr281@gmail.com" and "System" in the "Recorded By"
# column (G) wherever both names are present together.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 1) { $lastRow = 1 }

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
